# Update SituacionCalle StructureDefinition spreadsheet (0.2.1 -> refreshed
# publication): title/description text gets a space, the Context element
# reference is simplified, and the Date metadata is refreshed. "No es
# version nueva, solo actualizacion".

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Title (row 5) and Description (row 12) change from "SituacionCalle" to
# "Situacion Calle" (note the added space).
$meta.Range("B5").Value = "Situacion Calle"
$meta.Range("B12").Value = "Situacion Calle"

# Date metadata refreshed to the new publication timestamp.
$meta.Range("B8").Value = "2024-07-15T11:25:06-04:00"

# Context element simplified from the fully qualified URL to just "Address".
$meta.Range("B21").Value = "element:Address"

# Elements sheet: the root Extension row's Short/Definition columns (L2/M2)
# and the Extension.value[x] row's Short column (L6) get the same
# "SituacionCalle" -> "Situacion Calle" fix.
$elements.Range("L2").Value = "Situacion Calle"
$elements.Range("M2").Value = "Situacion Calle"
$elements.Range("L6").Value = "Situacion Calle"
